$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 180, pushing the existing weekly "Sandia" price rows
# (previously 180:213) down to 181:214, matching the inserted weekly entry
# pattern described in the commit ("Fruta / hortaliza, semanal").
$ws.Rows("180:180").Insert()

# Populate the newly inserted row 180 with this week's price data.
$ws.Range("A180").Value = 9
$ws.Range("B180").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C180").Value = "Metropolitana"
$ws.Range("D180").Value = 44504
$ws.Range("E180").Value = 13
$ws.Range("F180").Value = 100112028
$ws.Range("G180").Value = "Sandia"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 210
$ws.Range("K180").Value = 800
$ws.Range("L180").Value = 1000
$ws.Range("M180").Value = 900
$ws.Range("N180").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O180").Value = "Región de Arica y Parinacota"
$ws.Range("P180").Value = 900
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"
